$wb = $excel.ActiveWorkbook

function Set-AppendRow($ws, $r, $a, $b, $c, $d, $e, $f, $g, $h, $i) {
    $ws.Cells.Item($r, 1).Value = [double]$a
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = [double]$g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
}

$ws = $wb.Worksheets.Item(1)
Set-AppendRow $ws 159 "45945.49385416666" "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x00,0xB8" "0xf" 380 "7.598631275147109e+23" 196 15
Set-AppendRow $ws 160 "45946.49204861111" "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x00,0xB8" "0xf" 380 "7.598631275147109e+23" 192 15

$ws = $wb.Worksheets.Item(2)
Set-AppendRow $ws 159 "45945.49385416666" "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x00,0xBC" "0xe" 400 "5.68432987514711e+23" 196 14
Set-AppendRow $ws 160 "45946.49204861111" "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x00,0xB8" "0xe" 400 "5.68432987514711e+23" 196 14

$ws = $wb.Worksheets.Item(3)
Set-AppendRow $ws 159 "45945.49385416666" "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x50" "0x3" 110 "5.68631262647114e+23" 80 3
Set-AppendRow $ws 160 "45946.49204861111" "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x4F" "0x3" 110 "5.68631262647114e+23" 79 3

$ws = $wb.Worksheets.Item(4)
Set-AppendRow $ws 159 "45945.49385416666" "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x4C" "0x3" 110 "9.85046333984776e+23" 76 3
Set-AppendRow $ws 160 "45946.49204861111" "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x4C" "0x3" 110 "9.85046333984776e+23" 76 3
